$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 in I1 and IF in J1, matching the style of existing headers
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-15: I column is always 1, J column mirrors H column value
for ($r = 2; $r -le 15; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
